# Weekly update: insert two new "Plátano" price rows for Terminal La Palmera
# de La Serena (new reporting date 44449) above the existing historical
# series, pushing the rest of the series down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 292, shifting the
# existing data (old rows 292:326) down to 294:328.
$ws.Rows("292:293").Insert()

# --- New row 292: Pintón, 80 kg, date 44449 -------------------------------
$ws.Range("A292").Value = 8
$ws.Range("B292").Value = "Terminal La Palmera de La Serena"
$ws.Range("C292").Value = "Coquimbo"
$ws.Range("D292").Value = 44449
$ws.Range("E292").Value = 4
$ws.Range("F292").Value = "Fruta"
$ws.Range("G292").Value = 100108
$ws.Range("H292").Value = "Tropicales y subtropicales"
$ws.Range("I292").Value = 100108006
$ws.Range("J292").Value = "Plátano"
$ws.Range("K292").Value = "Sin especificar"
$ws.Range("L292").Value = "Pintón"
$ws.Range("M292").Value = 80
$ws.Range("N292").Value = 20000
$ws.Range("O292").Value = 20000
$ws.Range("P292").Value = 20000
$ws.Range("Q292").Value = "$/caja 20 kilos"
$ws.Range("R292").Value = "Ecuador"
$ws.Range("S292").Value = 1000
$ws.Range("T292").Value = 20

# --- New row 293: Primera Maduro, 120 kg, date 44449 ----------------------
$ws.Range("A293").Value = 8
$ws.Range("B293").Value = "Terminal La Palmera de La Serena"
$ws.Range("C293").Value = "Coquimbo"
$ws.Range("D293").Value = 44449
$ws.Range("E293").Value = 4
$ws.Range("F293").Value = "Fruta"
$ws.Range("G293").Value = 100108
$ws.Range("H293").Value = "Tropicales y subtropicales"
$ws.Range("I293").Value = 100108006
$ws.Range("J293").Value = "Plátano"
$ws.Range("K293").Value = "Sin especificar"
$ws.Range("L293").Value = "Primera Maduro"
$ws.Range("M293").Value = 120
$ws.Range("N293").Value = 22000
$ws.Range("O293").Value = 22000
$ws.Range("P293").Value = 22000
$ws.Range("Q293").Value = "$/caja 20 kilos"
$ws.Range("R293").Value = "Ecuador"
$ws.Range("S293").Value = 1100
$ws.Range("T293").Value = 20
